# Generate Report for Handoff
# - Drop the "faa2eff6-2430-4173-810c-d843ff92cbe2" item (now fully handed
#   off/in-sync) from every sheet, leaving only the "89f49ef1-..." item.
# - Refresh the status + timestamp of the remaining item to reflect the
#   new handoff.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-03-22 16:54:28"
$wsOverview.Range("A3").Hyperlinks.Delete()
$wsOverview.Rows.Item(3).Delete()

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-22 16:54:24"
$wsZhCn.Range("A3").Hyperlinks.Delete()
$wsZhCn.Range("D3").Hyperlinks.Delete()
$wsZhCn.Range("F3").Hyperlinks.Delete()
$wsZhCn.Range("G3").Hyperlinks.Delete()
$wsZhCn.Rows.Item(3).Delete()

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-22 16:54:28"
$wsDeDe.Range("A3").Hyperlinks.Delete()
$wsDeDe.Range("D3").Hyperlinks.Delete()
$wsDeDe.Range("F3").Hyperlinks.Delete()
$wsDeDe.Range("G3").Hyperlinks.Delete()
$wsDeDe.Rows.Item(3).Delete()
